# This edit cyclically re-shuffles the observation records currently stored
# in rows 2-7 of the active sheet ("artfynd" export). The set of rows/columns
# involved does not change; only the content of columns A, B, D, E, F, G, H,
# Q, R and S gets redistributed among rows 2-7 as follows:
#   new row 2 <- old row 6
#   new row 3 <- old row 7
#   new row 4 <- old row 2
#   new row 5 <- old row 4
#   new row 6 <- old row 5
#   new row 7 <- old row 3
# All other columns (C, I, J, K, P, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG,
# AT, AW, AX, AY, ...) are identical across these rows already, so they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "S")

# Mapping of destination row -> source row (values are read from the source
# row as it existed BEFORE this script runs).
$rowMap = @{
    2 = 6
    3 = 7
    4 = 2
    5 = 4
    6 = 5
    7 = 3
}

# Step 1: snapshot the current values of every involved cell (rows 2-7)
# before any writes happen, using Value2 so text/number round-trip cleanly.
$old = @{}
foreach ($r in 2..7) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $old[$addr] = $ws.Range($addr).Value2
    }
}

# Step 2: write the re-shuffled values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $destAddr = "$c$destRow"
        $srcAddr = "$c$srcRow"
        $ws.Range($destAddr).Value = $old[$srcAddr]
    }
}

Write-Output "Re-shuffled rows 2-7 (A,B,D,E,F,G,H,Q,R,S) per mapping $rowMap"
